$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-09-03 Tuesday"

# Update each arithmetic-problem cell in the table (row-major order)
$t = $d.Tables.Item(1)
$values = @(
    "15+16=",
    "86+9=",
    "63+28=",
    "62-4=",
    "16+7=",
    "61-5=",
    "31-25=",
    "38+48=",
    "37+28=",
    "91-45=",
    "5+28=",
    "44+7=",
    "16+28=",
    "94-85=",
    "83-66=",
    "14-9=",
    "93-87=",
    "19+26=",
    "85+6=",
    "69+13=",
    "90-9=",
    "84-15=",
    "12+79=",
    "46-18=",
    "93-49=",
    "93-27=",
    "73-14=",
    "5+36=",
    "40-35=",
    "96-69=",
    "86+6=",
    "61-55=",
    "45-39=",
    "93-58=",
    "76-7=",
    "29+24=",
    "23+8=",
    "53-26=",
    "34+17=",
    "76-28=",
    "92-25=",
    "54-19=",
    "7+48=",
    "62-48=",
    "32-14=",
    "34-16=",
    "72-29=",
    "95-79=",
    "82-69=",
    "79+4=",
    "69+8=",
    "27+69=",
    "49+32=",
    "71-22=",
    "63+29=",
    "91-87=",
    "52-29=",
    "78+8=",
    "35+49=",
    "18+27=",
    "24+18=",
    "91-62=",
    "42-34=",
    "84+7=",
    "4+88=",
    "98-39=",
    "57+39=",
    "62-57=",
    "75-7=",
    "13+79=",
    "83-57=",
    "38+34=",
    "5+48=",
    "54+7=",
    "61-29=",
    "23-14=",
    "96-39=",
    "44-8=",
    "48+17=",
    "14+47=",
    "43-15=",
    "19+28=",
    "80-71=",
    "77-8=",
    "54-28=",
    "14-9=",
    "45+29=",
    "73-49=",
    "9+76=",
    "17+35=",
    "91-25=",
    "63-59=",
    "70-43=",
    "93-65=",
    "33+9=",
    "15+76=",
    "63-55=",
    "70-33=",
    "67-8=",
    "26+67="
)

$numCols = $t.Columns.Count
$numRows = $t.Rows.Count
$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Output "Updated $idx cells; date set to $($d.Paragraphs.Item(1).Range.Text)"